$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 63÷2= -> 98÷9=
$t.Cell(1, 1).Range.Text = "98÷9="
# 87÷8= -> 72÷7=
$t.Cell(1, 2).Range.Text = "72÷7="
# 43÷5= -> 59÷2=
$t.Cell(1, 3).Range.Text = "59÷2="
# 60÷6= -> 33÷4=
$t.Cell(1, 4).Range.Text = "33÷4="
# 62÷8= -> 82÷8=
$t.Cell(1, 5).Range.Text = "82÷8="
# 60÷5= -> 41÷5=
$t.Cell(5, 1).Range.Text = "41÷5="
# 96÷5= -> 90÷6=
$t.Cell(5, 2).Range.Text = "90÷6="
# 68÷2= -> 99÷5=
$t.Cell(5, 3).Range.Text = "99÷5="
# 81÷6= -> 74÷7=
$t.Cell(5, 4).Range.Text = "74÷7="
# 61÷8= -> 17÷2=
$t.Cell(5, 5).Range.Text = "17÷2="
# 90÷6= -> 78÷5=
$t.Cell(9, 1).Range.Text = "78÷5="
# 67÷7= -> 70÷8=
$t.Cell(9, 2).Range.Text = "70÷8="
# 32÷9= -> 44÷2=
$t.Cell(9, 3).Range.Text = "44÷2="
# 62÷8= -> 80÷6=
$t.Cell(9, 4).Range.Text = "80÷6="
# 75÷8= -> 18÷2=
$t.Cell(9, 5).Range.Text = "18÷2="
# 14÷5= -> 57÷4=
$t.Cell(13, 1).Range.Text = "57÷4="
# 43÷8= -> 54÷7=
$t.Cell(13, 2).Range.Text = "54÷7="
# 53÷6= -> 13÷9=
$t.Cell(13, 3).Range.Text = "13÷9="
# 52÷2= -> 96÷3=
$t.Cell(13, 4).Range.Text = "96÷3="
# 78÷7= -> 25÷2=
$t.Cell(13, 5).Range.Text = "25÷2="
# 47÷8= -> 59÷8=
$t.Cell(17, 1).Range.Text = "59÷8="
# 89÷3= -> 65÷7=
$t.Cell(17, 2).Range.Text = "65÷7="
# 35÷3= -> 68÷2=
$t.Cell(17, 3).Range.Text = "68÷2="
# 81÷6= -> 43÷6=
$t.Cell(17, 4).Range.Text = "43÷6="
# 76÷4= -> 95÷7=
$t.Cell(17, 5).Range.Text = "95÷7="
